$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.495.16'
$ws.Range("E2").Value = '  +5.86%  '

$ws.Range("D3").Value = '2.591.34'
$ws.Range("E3").Value = '  +7.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.42'
$ws.Range("E5").Value = '  +3.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.07'
$ws.Range("E6").Value = '  +1.57%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -3.95%  '

$ws.Range("D9").Value = '2.621.25'
$ws.Range("E9").Value = '  +7.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.43'
$ws.Range("E10").Value = '  +5.10%  '

$ws.Range("E11").Value = '  +3.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("D14").Value = '3.050.34'
$ws.Range("E14").Value = '  +7.53%  '

$ws.Range("D15").Value = '60.474.68'
$ws.Range("E15").Value = '  +5.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.60'
$ws.Range("E16").Value = '  +4.84%  '

$ws.Range("E17").Value = '  +4.62%  '

$ws.Range("D18").Value = '2.614.11'
$ws.Range("E18").Value = '  +7.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.75'
$ws.Range("E19").Value = '  +2.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.13'
$ws.Range("E20").Value = '  +5.66%  '

$ws.Range("E21").Value = '  +4.06%  '

$ws.Range("E22").Value = '  +3.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.08'
$ws.Range("E24").Value = '  +3.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.422'
$ws.Range("E25").Value = '  +4.67%  '

$ws.Range("D26").Value = '2.722.75'
$ws.Range("E26").Value = '  +7.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.165'
$ws.Range("E27").Value = '  +2.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.991'
$ws.Range("E28").Value = '  -0.68%  '

$ws.Range("D29").Value = '0.0₃0855'
$ws.Range("E29").Value = '  +8.79%  '

$ws.Range("E30").Value = '  +3.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.43'
$ws.Range("E32").Value = '  +4.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.06'
$ws.Range("E33").Value = '  +3.01%  '

$ws.Range("E34").Value = '  +3.04%  '

$ws.Range("E35").Value = '  +7.83%  '

$ws.Range("E36").Value = '  +6.21%  '

$ws.Range("E37").Value = '  +3.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '308.26'
$ws.Range("E38").Value = '  +7.82%  '

$ws.Range("E39").Value = '  +7.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.844'
$ws.Range("E40").Value = '  +3.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.76'
$ws.Range("E41").Value = '  +6.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.833'
$ws.Range("E42").Value = '  +27.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.42'
$ws.Range("E43").Value = '  +4.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.626'
$ws.Range("E44").Value = '  +4.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0569'
$ws.Range("E45").Value = '  +7.13%  '

$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.995'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.78'
$ws.Range("E48").Value = '  +12.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.87'
$ws.Range("E49").Value = '  +6.66%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0235'
$ws.Range("E50").Value = '  +3.09%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.041.74'
$ws.Range("E51").Value = '  +6.71%  '

